$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "64.545.70"
Set-TextValue $ws "E2" "  -4.16%  "
Set-TextValue $ws "D3" "3.646.44"
Set-TextValue $ws "E3" "  +1.92%  "
Set-TextValue $ws "D4" "1.01"
Set-TextValue $ws "E4" "  +0.89%  "
Set-TextValue $ws "D5" "404.52"
Set-TextValue $ws "E5" "  -3.18%  "
Set-TextValue $ws "D6" "131.13"
Set-TextValue $ws "E6" "  +0.44%  "
Set-TextValue $ws "D7" "3.630.09"
Set-TextValue $ws "E7" "  -3.81%  "
Set-TextValue $ws "D8" "0.619"
Set-TextValue $ws "E8" "  -5.12%  "
Set-TextValue $ws "E9" "  +0.27%  "
Set-TextValue $ws "D10" "0.719"
Set-TextValue $ws "E10" "  -8.48%  "
Set-TextValue $ws "D11" "0.159"
Set-TextValue $ws "E11" "  -13.09%  "
Set-TextValue $ws "D12" "0.0000295"
Set-TextValue $ws "E12" "  -14.95%  "
Set-TextValue $ws "D13" "41.64"
Set-TextValue $ws "E13" "  -2.96%  "
Set-TextValue $ws "D14" "9.84"
Set-TextValue $ws "E14" "  -1.50%  "
Set-TextValue $ws "D15" "4.236.49"
Set-TextValue $ws "E15" "  +2.47%  "
Set-TextValue $ws "E16" "  -1.48%  "
Set-TextValue $ws "B17" "Uniswap"
Set-TextValue $ws "C17" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws "D17" "13.47"
Set-TextValue $ws "E17" "  +7.62%  "
Set-TextValue $ws "B18" "WrappedEther"
Set-TextValue $ws "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D18" "3.600.37"
Set-TextValue $ws "E18" "  +0.53%  "
Set-TextValue $ws "B19" "Chainlink"
Set-TextValue $ws "C19" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D19" "19.75"
Set-TextValue $ws "E19" "  -2.52%  "
Set-TextValue $ws "D20" "1.07"
Set-TextValue $ws "E20" "  -5.17%  "
Set-TextValue $ws "D21" "64.843.38"
Set-TextValue $ws "E21" "  -3.28%  "
Set-TextValue $ws "D22" "417.11"
Set-TextValue $ws "E22" "  -9.66%  "
Set-TextValue $ws "D23" "15.31"
Set-TextValue $ws "E23" "  +17.36%  "
Set-TextValue $ws "D24" "85.31"
Set-TextValue $ws "E24" "  -5.53%  "
Set-TextValue $ws "D25" "2.97"
Set-TextValue $ws "E25" "  -6.77%  "
Set-TextValue $ws "D26" "35.58"
Set-TextValue $ws "E26" "  +1.46%  "
Set-TextValue $ws "D27" "3.15"
Set-TextValue $ws "E27" "  -6.71%  "
Set-TextValue $ws "D28" "9.32"
Set-TextValue $ws "E28" "  -6.14%  "
Set-TextValue $ws "D29" "5.11"
Set-TextValue $ws "E29" "  +5.27%  "
Set-TextValue $ws "E30" "  +0.11%  "
Set-TextValue $ws "D31" "2.69"
Set-TextValue $ws "E31" "  -3.84%  "
Set-TextValue $ws "D32" "0.117"
Set-TextValue $ws "E32" "  -0.28%  "
Set-TextValue $ws "D33" "6.90"
Set-TextValue $ws "E33" "  -5.72%  "
Set-TextValue $ws "E34" "  +1.39%  "
Set-TextValue $ws "D35" "40.23"
Set-TextValue $ws "E35" "  +2.15%  "
Set-TextValue $ws "D36" "55.77"
Set-TextValue $ws "E36" "  -1.43%  "
Set-TextValue $ws "E37" "  -0.09%  "
Set-TextValue $ws "D38" "0.0459"
Set-TextValue $ws "E38" "  -7.33%  "
Set-TextValue $ws "D39" "2.88"
Set-TextValue $ws "E39" "  +25.88%  "
Set-TextValue $ws "D40" "0.998"
Set-TextValue $ws "E40" "  +0.31%  "
Set-TextValue $ws "E41" "  -6.44%  "
Set-TextValue $ws "D42" "27.05"
Set-TextValue $ws "E42" "  +25.39%  "
Set-TextValue $ws "D43" "3.30"
Set-TextValue $ws "E43" "  +1.94%  "
Set-TextValue $ws "B44" "NEARProtocol"
Set-TextValue $ws "C44" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D44" "4.34"
Set-TextValue $ws "E44" "  -0.40%  "
Set-TextValue $ws "B45" "Monero"
Set-TextValue $ws "C45" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D45" "141.85"
Set-TextValue $ws "E45" "  -4.53%  "
Set-TextValue $ws "B46" "ARBITRUM"
Set-TextValue $ws "C46" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D46" "2.05"
Set-TextValue $ws "E46" "  +3.58%  "
Set-TextValue $ws "B47" "PEPE"
Set-TextValue $ws "C47" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws "D47" "0.0₃0614"
Set-TextValue $ws "E47" "  -22.79%  "
Set-TextValue $ws "B48" "ApeXProtocol"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D48" "3.03"
Set-TextValue $ws "E48" "  +16.91%  "
Set-TextValue $ws "D49" "2.78"
Set-TextValue $ws "E49" "  -6.95%  "
Set-TextValue $ws "E50" "  -8.12%  "
Set-TextValue $ws "D51" "0.288"
Set-TextValue $ws "E51" "  -6.47%  "
